# aggiornamento fino a 28 luglio
# Appends new daily rows (302-328) to the Carpi report sheet: date serial,
# "nuovi pos." count, the rolling 7-day sum, and the rolling 7-day sum per
# 100k inhabitants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 302

$data = @(
  @(44376, 0, 1, 1.381807127361163),
  @(44377, 0, 0, 0),
  @(44378, 3, 3, 4.145421382083488),
  @(44379, 0, 3, 4.145421382083488),
  @(44380, 1, 4, 5.527228509444652),
  @(44381, 0, 4, 5.527228509444652),
  @(44382, 1, 5, 6.909035636805815),
  @(44383, 0, 5, 6.909035636805815),
  @(44384, 0, 5, 6.909035636805815),
  @(44385, 0, 2, 2.763614254722326),
  @(44386, 0, 2, 2.763614254722326),
  @(44387, 0, 1, 1.381807127361163),
  @(44388, 0, 1, 1.381807127361163),
  @(44389, 2, 2, 2.763614254722326),
  @(44390, 0, 2, 2.763614254722326),
  @(44391, 1, 3, 4.145421382083488),
  @(44392, 1, 4, 5.527228509444652),
  @(44393, 0, 4, 5.527228509444652),
  @(44394, 3, 7, 9.672649891528142),
  @(44395, 0, 7, 9.672649891528142),
  @(44396, 0, 5, 6.909035636805815),
  @(44397, 0, 5, 6.909035636805815),
  @(44398, 0, 4, 5.527228509444652),
  @(44399, 5, 8, 11.0544570188893),
  @(44400, 6, 14, 19.34529978305628),
  @(44401, 18, 29, 40.07240669347372),
  @(44402, 19, 48, 66.32674211333581)
)

$endRow = $startRow + $data.Count - 1

# Column A carries the same date style (s="2": centered, bordered, custom
# date numFmt) as every other row in the sheet - copy it down from the last
# existing row instead of re-declaring it, so no duplicate style gets
# created.
$ws.Cells.Item($startRow - 1, 1).Copy()
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $row = $data[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
}
